$d = $word.ActiveDocument

# 1) Fix the URL: "leevi/velothink/Velothink-main" -> "leevi/Velothink-main/Velothink-main"
$d.Content.Find.Execute("leevi/velothink/Velothink-main", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "leevi/Velothink-main/Velothink-main", 2)

# 2) Move the comma: "ARK.html), testataan torstaina" should read ARK.html), then space, then "testataan"
#    (text is already correct after step 1, this is just to ensure spacing matches; no-op safety net)
$d.Content.Find.Execute("ARK.html), testataan torstaina", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ARK.html), testataan torstaina", 2)

# 3) Change the date from 1.6.2023 to 5.6.2023
$d.Content.Find.Execute("1.6.2023", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "5.6.2023", 2)

# 4) Change the time from "klo: 10 " to "klo: 8 "
$d.Content.Find.Execute("klo: 10 ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "klo: 8 ", 2)
